$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.784.53"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "3.952.86"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "605.20"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "173.27"
$ws.Range("E6").Value = "  +12.40%  "
$ws.Range("D7").Value = "0.683"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.789"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").Value = "0.182"
$ws.Range("E10").Value = "  +6.81%  "
$ws.Range("D11").Value = "56.29"
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("D12").Value = "0.0000328"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "11.61"
$ws.Range("E13").Value = "  +5.14%  "
$ws.Range("D14").Value = "4.573.31"
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").Value = "3.961.31"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "21.55"
$ws.Range("E16").Value = "  +3.89%  "
$ws.Range("D17").Value = "14.08"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "72.648.04"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "445.54"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "4.80"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").Value = "95.75"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("E24").Value = "  -5.35%  "
$ws.Range("D25").Value = "14.20"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "11.23"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "5.90"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "10.42"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").Value = "35.94"
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").Value = "13.96"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "49.71"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("D35").Value = "0.0₃0994"
$ws.Range("E35").Value = "  +13.57%  "
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("E37").Value = "  -8.48%  "
$ws.Range("D38").Value = "0.428"
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("D39").Value = "3.44"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "3.25"
$ws.Range("E43").Value = "  +44.13%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0478"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D45").Value = "10.61"
$ws.Range("E45").Value = "  -6.43%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "2.64"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("D48").Value = "3.40"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000286"
$ws.Range("E49").Value = "  +5.89%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "2.77"
$ws.Range("E50").Value = "  -18.75%  "
$ws.Range("D51").Value = "2.829.68"
$ws.Range("E51").Value = "  +1.42%  "
